$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated price (D) and volume-change (E) figures, and the
# THORChain/PaxDollar row swap (rows 49-50), per the latest scrape.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.907.09'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.61%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.840.43'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.49%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.90'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.92%  '
$ws.Range("E9").Value = '  +1.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0688'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0985'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.107.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("E13").Value = '  +3.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.838.72'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.674'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.58%  '
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '34.887.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '240.64'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.62%  '
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.49'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.80'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("E28").Value = '  +2.48%  '
$ws.Range("E29").Value = '  -5.57%  '
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("E32").Value = '  -6.10%  '
$ws.Range("E33").Value = '  -1.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.92'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.54%  '
$ws.Range("E35").Value = '  +8.49%  '
$ws.Range("E36").Value = '  +11.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.697'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.99%  '
$ws.Range("E38").Value = '  +6.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '91.12'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.343.30'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.53%  '
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.85'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("E44").Value = '  -2.88%  '
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("E47").Value = '  +2.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.021.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.50%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.42'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +20.15%  '
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("E51").Value = '  +2.16%  '
